$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.829.71"
$ws.Range("E2").Value = "  +5.73%  "

$ws.Range("D3").Value = "3.415.33"
$ws.Range("E3").Value = "  +11.48%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.90"
$ws.Range("E5").Value = "  +2.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.23"
$ws.Range("E6").Value = "  +8.76%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "3.407.77"
$ws.Range("E8").Value = "  +11.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  +4.85%  "

$ws.Range("E10").Value = "  +3.46%  "

$ws.Range("E11").Value = "  +5.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.489"
$ws.Range("E12").Value = "  +4.31%  "

$ws.Range("E13").Value = "  +5.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.47"
$ws.Range("E14").Value = "  +7.95%  "

$ws.Range("D15").Value = "3.970.37"
$ws.Range("E15").Value = "  +11.04%  "

$ws.Range("D16").Value = "69.679.10"
$ws.Range("E16").Value = "  +5.49%  "

$ws.Range("E17").Value = "  +1.43%  "

$ws.Range("D18").Value = "3.402.30"
$ws.Range("E18").Value = "  +10.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.43"
$ws.Range("E19").Value = "  +7.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.11"
$ws.Range("E20").Value = "  +2.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "504.23"
$ws.Range("E21").Value = "  +3.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.39"
$ws.Range("E22").Value = "  +9.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("E23").Value = "  +6.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.55"
$ws.Range("E24").Value = "  +4.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.31"
$ws.Range("E25").Value = "  +5.87%  "

$ws.Range("E26").Value = "  +8.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.81"
$ws.Range("E27").Value = "  +6.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.27"
$ws.Range("E29").Value = "  +5.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.51"
$ws.Range("E30").Value = "  +11.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.70"
$ws.Range("E31").Value = "  +3.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.74"
$ws.Range("E32").Value = "  +7.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0000103"
$ws.Range("E33").Value = "  +14.59%  "

$ws.Range("E34").Value = "  +4.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.06"
$ws.Range("E36").Value = "  +8.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +6.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.22"
$ws.Range("E38").Value = "  +6.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.330"
$ws.Range("E39").Value = "  +10.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.11"
$ws.Range("E40").Value = "  +7.73%  "

$ws.Range("E41").Value = "  +6.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "50.09"
$ws.Range("E42").Value = "  +2.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.67"
$ws.Range("E43").Value = "  +5.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  +12.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "413.11"
$ws.Range("E45").Value = "  +14.01%  "

$ws.Range("D46").Value = "2.936.44"
$ws.Range("E46").Value = "  +5.43%  "

$ws.Range("E47").Value = "  +4.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.72"
$ws.Range("E48").Value = "  +14.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "135.30"
$ws.Range("E49").Value = "  +0.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.44"
$ws.Range("E51").Value = "  +14.13%  "

